$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.696.46'
$ws.Range('E2').Value = '  -3.75%  '

# Row 3
$ws.Range('D3').Value = '3.880.89'
$ws.Range('E3').Value = '  -3.86%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
$ws.Range('D5').Value = '528.23'
$ws.Range('E5').Value = '  +1.83%  '

# Row 6
$ws.Range('D6').Value = '143.98'
$ws.Range('E6').Value = '  -1.76%  '

# Row 7
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.872.51'
$ws.Range('E7').Value = '  -3.76%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.671'
$ws.Range('E8').Value = '  -8.81%  '

# Row 9
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.13%  '

# Row 10
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = '0.720'
$ws.Range('E10').Value = '  -5.25%  '

# Row 11
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.163'
$ws.Range('E11').Value = '  -5.97%  '

# Row 12
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').Value = '52.59'
$ws.Range('E12').Value = '  +11.91%  '

# Row 13
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '0.0000308'
$ws.Range('E13').Value = '  -4.97%  '

# Row 14
$ws.Range('D14').Value = '4.511.41'
$ws.Range('E14').Value = '  -3.69%  '

# Row 15
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '10.27'
$ws.Range('E15').Value = '  -5.84%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.883.20'
$ws.Range('E16').Value = '  -3.50%  '

# Row 17
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '13.56'
$ws.Range('E17').Value = '  -3.77%  '

# Row 18
$ws.Range('E18').Value = '  -1.22%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '19.89'
$ws.Range('E19').Value = '  -5.77%  '

# Row 20
$ws.Range('B20').Value = 'Polygon'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D20').Value = '1.15'
$ws.Range('E20').Value = '  -5.21%  '

# Row 21
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '69.719.63'
$ws.Range('E21').Value = '  -3.41%  '

# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '419.68'
$ws.Range('E22').Value = '  -4.98%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '94.50'
$ws.Range('E23').Value = '  -9.47%  '

# Row 24
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '3.43'
$ws.Range('E24').Value = '  -4.63%  '

# Row 25
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '4.07'
$ws.Range('E25').Value = '  +1.91%  '

# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '13.87'
$ws.Range('E26').Value = '  -4.93%  '

# Row 27
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '10.95'
$ws.Range('E27').Value = '  -4.52%  '

# Row 28
$ws.Range('D28').Value = '5.84'
$ws.Range('E28').Value = '  +0.62%  '

# Row 29
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = '10.29'
$ws.Range('E29').Value = '  -6.89%  '

# Row 30
$ws.Range('D30').Value = '3.53'
$ws.Range('E30').Value = '  +12.72%  '

# Row 31
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '35.56'
$ws.Range('E31').Value = '  -6.03%  '

# Row 32
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '7.19'
$ws.Range('E32').Value = '  +5.34%  '

# Row 33
$ws.Range('D33').Value = '676.96'
$ws.Range('E33').Value = '  +0.67%  '

# Row 34
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').Value = '12.94'
$ws.Range('E34').Value = '  -5.14%  '

# Row 35
$ws.Range('D35').Value = '47.10'
$ws.Range('E35').Value = '  +11.00%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.125'
$ws.Range('E36').Value = '  -2.36%  '

# Row 37
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = '64.09'
$ws.Range('E37').Value = '  -4.74%  '

# Row 38
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').Value = '0.419'
$ws.Range('E38').Value = '  -2.77%  '

# Row 39
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  -4.18%  '

# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.16%  '

# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.145'
$ws.Range('E41').Value = '  -3.47%  '

# Row 42
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '0.0₃0790'
$ws.Range('E42').Value = '  -7.71%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.09%  '

# Row 44
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '3.22'
$ws.Range('E44').Value = '  +0.30%  '

# Row 45
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0469'
$ws.Range('E45').Value = '  -5.17%  '

# Row 46
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.145'
$ws.Range('E46').Value = '  -9.13%  '

# Row 47
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = '2.61'
$ws.Range('E47').Value = '  -3.45%  '

# Row 48
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '9.44'
$ws.Range('E48').Value = '  +4.59%  '

# Row 49
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '3.31'
$ws.Range('E49').Value = '  -4.58%  '

# Row 50
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  -5.87%  '

# Row 51
$ws.Range('D51').Value = '143.97'
$ws.Range('E51').Value = '  +0.97%  '
